$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D (new quarterly period),
# shifting existing quarterly data from D:K to F:M.
$ws.Range("D5:E102").EntireColumn.Insert()

# Carry over the number formatting/styles from the old D:E columns
# (now shifted to F:G) into the newly inserted D:E columns.
$ws.Range("F5:G102").Copy()
$ws.Range("D5:E102").PasteSpecial(-4122)

# A few rows have no quarterly data at all (section headers / blank
# separator rows); remove the blank D:E stubs the paste above created
# there so they stay empty just like in the source rows.
$ws.Range("D5:E6").Clear()
$ws.Range("D36:E37").Clear()
$ws.Range("D78:E79").Clear()

# Populate the new columns with the latest reported quarter data,
# and apply a couple of restated figures in the Capital Expenditures row.
$cellValues = @(
    @("D7", 43465),
    @("E7", 43373),
    @("D8", 143400),
    @("E8", 153500),
    @("D9", 500),
    @("E9", 400),
    @("D10", 142900),
    @("E10", 153100),
    @("D12", "NA"),
    @("E12", "NA"),
    @("D13", 0),
    @("E13", 0),
    @("D14", -47100),
    @("E14", 0),
    @("D15", 16600),
    @("E15", 37000),
    @("D17", -6400),
    @("E17", 96300),
    @("D18", 149800),
    @("E18", 57200),
    @("D20", 2800),
    @("E20", 800),
    @("D21", 187600),
    @("E21", 95000),
    @("D22", 13500),
    @("E22", 11600),
    @("D23", 139000),
    @("E23", 46400),
    @("D24", -300),
    @("E24", 100),
    @("D25", 0),
    @("E25", 0),
    @("D26", 139300),
    @("E26", 46300),
    @("D27", 139300),
    @("E27", 46300),
    @("D28", 0),
    @("E28", 0),
    @("D29", -500),
    @("E29", "NA"),
    @("D30", 0),
    @("E30", 0),
    @("D31", 0),
    @("E31", 0),
    @("D32", -2800),
    @("E32", -800),
    @("D33", 138800),
    @("E33", 46300),
    @("D34", 0),
    @("E34", 0),
    @("D35", 138800),
    @("E35", 46300),
    @("D38", 43465),
    @("E38", 43373),
    @("D41", 33300),
    @("E41", 339100),
    @("D42", 0),
    @("E42", 0),
    @("D43", 116500),
    @("E43", 131200),
    @("D44", 0),
    @("E44", 0),
    @("D45", 76400),
    @("E45", 19700),
    @("D46", 226200),
    @("E46", 490000),
    @("D47", 0),
    @("E47", 0),
    @("D48", 515400),
    @("E48", 522800),
    @("D49", 0),
    @("E49", 0),
    @("D50", 0),
    @("E50", 0),
    @("D51", 0),
    @("E51", 0),
    @("D52", 107200),
    @("E52", 89600),
    @("D53", 0),
    @("E53", 0),
    @("D54", 848900),
    @("E54", 1102300),
    @("D57", 82100),
    @("E57", 95500),
    @("D58", 0),
    @("E58", 224800),
    @("D59", 104200),
    @("E59", 126500),
    @("D60", 186300),
    @("E60", 446800),
    @("D61", 633500),
    @("E61", 759100),
    @("D62", 353800),
    @("E62", 356200),
    @("D63", 0),
    @("E63", 0),
    @("D64", 0),
    @("E64", 0),
    @("D65", 0),
    @("E65", 0),
    @("D66", 1173700),
    @("E66", 1562100),
    @("D68", 0),
    @("E68", 0),
    @("D69", 0),
    @("E69", 0),
    @("D70", 0),
    @("E70", 0),
    @("D71", 0),
    @("E71", 0),
    @("D72", -846300),
    @("E72", -985200),
    @("D73", 0),
    @("E73", 0),
    @("D74", 0),
    @("E74", 0),
    @("D75", 0),
    @("E75", 0),
    @("D76", -324800),
    @("E76", -459800),
    @("D77", 0),
    @("E77", 0),
    @("D80", 43465),
    @("E80", 43373),
    @("D81", 138800),
    @("E81", 46300),
    @("D83", 35000),
    @("E83", 37000),
    @("D84", 0),
    @("E84", 0),
    @("D85", 0),
    @("E85", 0),
    @("D86", 0),
    @("E86", 0),
    @("D87", 0),
    @("E87", 0),
    @("D88", 0),
    @("E88", 0),
    @("D89", 26900),
    @("E89", 179700),
    @("D91", 0),
    @("E91", -200),
    @("F91", 4500),
    @("H91", 0),
    @("I91", -100),
    @("J91", 0),
    @("D92", 0),
    @("E92", 0),
    @("D93", 0),
    @("E93", 0),
    @("D94", -20700),
    @("E94", 32000),
    @("D96", 0),
    @("E96", 0),
    @("D97", 0),
    @("E97", 0),
    @("D98", 0),
    @("E98", 0),
    @("D99", 0),
    @("E99", 0),
    @("D100", -312000),
    @("E100", -2100),
    @("D101", 0),
    @("E101", 0),
    @("D102", -305800),
    @("E102", 209600)
)

foreach ($pair in $cellValues) {
    $ws.Range($pair[0]).Value = $pair[1]
}
